$d = $word.ActiveDocument

# 1) "qCaves" -> "Caves" (drop the stray leading "q" run/spell-check markers)
#    This paragraph is the unique one whose text is exactly "qCaves".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r") -eq "qCaves") {
        $p.Range.Find.Execute("qCaves", $true, $false, $false, $false, $false, $true, 1, $false, "Caves", 2)
        break
    }
}

# 2) Append shop prices to the "[In shop]" boost entries under "RPG Combat"
#    (By level -> Castle Approach / Castle Interior / Sundered Castle / Caves / Dev Wasteland).
#    Each entry is matched by its unique preceding label text within that paragraph,
#    since "[In shop]" itself is not unique across the whole document.
$shopFixes = @(
    @{ Label = "HP +10 ";      Price = 50 },
    @{ Label = "Strength +4 "; Price = 50 },
    @{ Label = "Magic +4 ";    Price = 50 },
    @{ Label = "HP +15 ";      Price = 75 },
    @{ Label = "Strength +4 "; Price = 50 },
    @{ Label = "HP +10 ";      Price = 50 },
    @{ Label = "MP +2 ";       Price = 50 },
    @{ Label = "Magic +4 ";    Price = 50 },
    @{ Label = "HP +10 ";      Price = 50 },
    @{ Label = "Strength +2 "; Price = 25 },
    @{ Label = "Magic +4 ";    Price = 50 },
    @{ Label = "HP +10 ";      Price = 50 },
    @{ Label = "MP +2 ";       Price = 60 },
    @{ Label = "Strength +5 "; Price = 60 },
    @{ Label = "Magic +5 ";    Price = 60 }
)

# Find the "RPG Combat" section header paragraph; the fixes start right after it.
$startIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd("`r") -eq "RPG Combat [In shop]") {
        $startIndex = $i + 1
        break
    }
}

$fixPos = 0
for ($i = $startIndex; $i -le $d.Paragraphs.Count -and $fixPos -lt $shopFixes.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text.TrimEnd("`r")
    $fix = $shopFixes[$fixPos]
    $expected = $fix.Label + "[In shop]"
    if ($text -eq $expected) {
        $newText = "[In shop: " + $fix.Price + "]"
        $p.Range.Find.Execute("[In shop]", $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
        $fixPos = $fixPos + 1
    }
}

Write-Host "qCaves fix done; shop price fixes applied: $fixPos of $($shopFixes.Count)"
